$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new "Numero" (A), new "Groupe" (D), new "Moyenne de l'etudiant" (E)
$data = @(
    @{Row=3; A=20170926; D="1-B"; E=8},
    @{Row=4; A=20170927; D="1-C"; E=11},
    @{Row=5; A=20170928; D="1-A"; E=6},
    @{Row=6; A=20170929; D="1-A"; E=6},
    @{Row=7; A=20170930; D="1-A"; E=14},
    @{Row=8; A=20170931; D="1-B"; E=19},
    @{Row=9; A=20170932; D="1-C"; E=14},
    @{Row=10; A=20170933; D="1-C"; E=11},
    @{Row=11; A=20170934; D="1-B"; E=15},
    @{Row=12; A=20170935; D="1-C"; E=10},
    @{Row=13; A=20170936; D="1-C"; E=17},
    @{Row=14; A=20170937; D="1-C"; E=18},
    @{Row=15; A=20170938; D="1-B"; E=5},
    @{Row=16; A=20170939; D="1-A"; E=8},
    @{Row=17; A=20170940; D="1-C"; E=16},
    @{Row=18; A=20170941; D="1-C"; E=11},
    @{Row=19; A=20170942; D="1-B"; E=17},
    @{Row=20; A=20170943; D="1-B"; E=19},
    @{Row=21; A=20170944; D="1-B"; E=15},
    @{Row=22; A=20170945; D="1-C"; E=7},
    @{Row=23; A=20170946; D="1-A"; E=8},
    @{Row=24; A=20170947; D="1-A"; E=13},
    @{Row=25; A=20170948; D="1-B"; E=5},
    @{Row=26; A=20170949; D="1-B"; E=15},
    @{Row=27; A=20170950; D="1-A"; E=12},
    @{Row=28; A=20170951; D="1-C"; E=10},
    @{Row=29; A=20170952; D="1-B"; E=17},
    @{Row=30; A=20170953; D="1-B"; E=13},
    @{Row=31; A=20170954; D="1-A"; E=17},
    @{Row=32; A=20170955; D="1-C"; E=5},
    @{Row=33; A=20170956; D="1-B"; E=9},
    @{Row=34; A=20170957; D="1-C"; E=17},
    @{Row=35; A=20170958; D="1-A"; E=19},
    @{Row=36; A=20170959; D="1-C"; E=13},
    @{Row=37; A=20170960; D="1-A"; E=13},
    @{Row=38; A=20170961; D="1-A"; E=9},
    @{Row=39; A=20170962; D="1-A"; E=13},
    @{Row=40; A=20170963; D="1-B"; E=16},
    @{Row=41; A=20170964; D="1-A"; E=7},
    @{Row=42; A=20170965; D="1-C"; E=6},
    @{Row=43; A=20170966; D="1-A"; E=16},
    @{Row=44; A=20170967; D="1-A"; E=8},
    @{Row=45; A=20170968; D="1-C"; E=8},
    @{Row=46; A=20170969; D="1-A"; E=20},
    @{Row=47; A=20170970; D="1-B"; E=10},
    @{Row=48; A=20170971; D="1-C"; E=16},
    @{Row=49; A=20170972; D="1-B"; E=11},
    @{Row=50; A=20170973; D="1-B"; E=13},
    @{Row=51; A=20170974; D="1-C"; E=14},
    @{Row=52; A=20170975; D="1-A"; E=20},
    @{Row=53; A=20170976; D="1-B"; E=13},
    @{Row=54; A=20170977; D="1-C"; E=16},
    @{Row=55; A=20170978; D="1-C"; E=6},
    @{Row=56; A=20170979; D="1-B"; E=5},
    @{Row=57; A=20170980; D="1-A"; E=5},
    @{Row=58; A=20170981; D="1-B"; E=20},
    @{Row=59; A=20170982; D="1-B"; E=7},
    @{Row=60; A=20170983; D="1-A"; E=17},
    @{Row=61; A=20170984; D="1-C"; E=19},
    @{Row=62; A=20170985; D="1-A"; E=20},
    @{Row=63; A=20170986; D="1-B"; E=16}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

Write-Output "Updated $($data.Count) rows"
